# Add Portuguese ("por") translation row to the machine_spec master-data
# sheet, mirroring the existing "eng" row (row 2) with translated id/name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "por"
$ws.Range("B3").Value = "RESIDENTE-1"
$ws.Range("C3").Value = "Máquina Virtual Residente"
$ws.Range("D3").Value = "Unkown"
$ws.Range("E3").Value = "Unknown"
$ws.Range("F3").Value = "RESIDENT-REG"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "Resident Virtual Machine"
$ws.Range("I3").Value = $true
$ws.Range("I3").NumberFormat = $ws.Range("I2").NumberFormat

$ws.Range("B3").WrapText = $true

$ws.Range("C4").Select()
